$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Live Data")

$ws.Cells.Item(2, 1).Value = "Bitcoin"
$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = 84016
$ws.Cells.Item(2, 4).Value = 1661676879344
$ws.Cells.Item(2, 5).Value = 69355645598
$ws.Cells.Item(2, 6).Value = -7.02592

$ws.Cells.Item(3, 1).Value = "Ethereum"
$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = 2080.51
$ws.Cells.Item(3, 4).Value = 250862158335
$ws.Cells.Item(3, 5).Value = 36605874923
$ws.Cells.Item(3, 6).Value = -8.74573

$ws.Cells.Item(4, 1).Value = "Tether"
$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = 0.999385
$ws.Cells.Item(4, 4).Value = 142420402555
$ws.Cells.Item(4, 5).Value = 114682656494
$ws.Cells.Item(4, 6).Value = -0.00698

$ws.Cells.Item(5, 1).Value = "XRP"
$ws.Cells.Item(5, 2).Value = "XRP"
$ws.Cells.Item(5, 3).Value = 2.4
$ws.Cells.Item(5, 4).Value = 138450797233
$ws.Cells.Item(5, 5).Value = 11278322557
$ws.Cells.Item(5, 6).Value = -7.86962

$ws.Cells.Item(6, 1).Value = "BNB"
$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = 571.1799999999999
$ws.Cells.Item(6, 4).Value = 83202775050
$ws.Cells.Item(6, 5).Value = 1499976633
$ws.Cells.Item(6, 6).Value = -4.05883

$ws.Cells.Item(7, 1).Value = "Solana"
$ws.Cells.Item(7, 2).Value = "SOL"
$ws.Cells.Item(7, 3).Value = 137.76
$ws.Cells.Item(7, 4).Value = 69895832790
$ws.Cells.Item(7, 5).Value = 10102969324
$ws.Cells.Item(7, 6).Value = -12.72925

$ws.Cells.Item(8, 1).Value = "USDC"
$ws.Cells.Item(8, 2).Value = "USDC"
$ws.Cells.Item(8, 3).Value = 0.999877
$ws.Cells.Item(8, 4).Value = 56285902585
$ws.Cells.Item(8, 5).Value = 17240121458
$ws.Cells.Item(8, 6).Value = -0.00427

$ws.Cells.Item(9, 1).Value = "Cardano"
$ws.Cells.Item(9, 2).Value = "ADA"
$ws.Cells.Item(9, 3).Value = 0.880883
$ws.Cells.Item(9, 4).Value = 31482021533
$ws.Cells.Item(9, 5).Value = 4748013946
$ws.Cells.Item(9, 6).Value = -9.19943

$ws.Cells.Item(10, 1).Value = "Dogecoin"
$ws.Cells.Item(10, 2).Value = "DOGE"
$ws.Cells.Item(10, 3).Value = 0.192464
$ws.Cells.Item(10, 4).Value = 28439433598
$ws.Cells.Item(10, 5).Value = 2778852678
$ws.Cells.Item(10, 6).Value = -10.92039

$ws.Cells.Item(11, 1).Value = "TRON"
$ws.Cells.Item(11, 2).Value = "TRX"
$ws.Cells.Item(11, 3).Value = 0.237222
$ws.Cells.Item(11, 4).Value = 20390178189
$ws.Cells.Item(11, 5).Value = 1510856163
$ws.Cells.Item(11, 6).Value = -1.82534

$ws.Cells.Item(12, 1).Value = "Lido Staked Ether"
$ws.Cells.Item(12, 2).Value = "STETH"
$ws.Cells.Item(12, 3).Value = 2081.02
$ws.Cells.Item(12, 4).Value = 19502472120
$ws.Cells.Item(12, 5).Value = 188780856
$ws.Cells.Item(12, 6).Value = -8.40056

$ws.Cells.Item(13, 1).Value = "Pi Network"
$ws.Cells.Item(13, 2).Value = "PI"
$ws.Cells.Item(13, 3).Value = 1.75
$ws.Cells.Item(13, 4).Value = 12134273796
$ws.Cells.Item(13, 5).Value = 815815324
$ws.Cells.Item(13, 6).Value = 1.23768

$ws.Cells.Item(14, 1).Value = "Wrapped Bitcoin"
$ws.Cells.Item(14, 2).Value = "WBTC"
$ws.Cells.Item(14, 3).Value = 83683
$ws.Cells.Item(14, 4).Value = 10801840781
$ws.Cells.Item(14, 5).Value = 808245158
$ws.Cells.Item(14, 6).Value = -7.03777

$ws.Cells.Item(15, 1).Value = "Hedera"
$ws.Cells.Item(15, 2).Value = "HBAR"
$ws.Cells.Item(15, 3).Value = 0.236245
$ws.Cells.Item(15, 4).Value = 9875467986
$ws.Cells.Item(15, 5).Value = 896388403
$ws.Cells.Item(15, 6).Value = -3.02431

$ws.Cells.Item(16, 1).Value = "LEO Token"
$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = 9.949999999999999
$ws.Cells.Item(16, 4).Value = 9192856882
$ws.Cells.Item(16, 5).Value = 9472624
$ws.Cells.Item(16, 6).Value = -0.24846

$ws.Cells.Item(17, 1).Value = "Wrapped stETH"
$ws.Cells.Item(17, 2).Value = "WSTETH"
$ws.Cells.Item(17, 3).Value = 2483.44
$ws.Cells.Item(17, 4).Value = 8816910376
$ws.Cells.Item(17, 5).Value = 83748391
$ws.Cells.Item(17, 6).Value = -8.757569999999999

$ws.Cells.Item(18, 1).Value = "Chainlink"
$ws.Cells.Item(18, 2).Value = "LINK"
$ws.Cells.Item(18, 3).Value = 13.78
$ws.Cells.Item(18, 4).Value = 8786978171
$ws.Cells.Item(18, 5).Value = 1097688210
$ws.Cells.Item(18, 6).Value = -14.76197

$ws.Cells.Item(19, 1).Value = "Stellar"
$ws.Cells.Item(19, 2).Value = "XLM"
$ws.Cells.Item(19, 3).Value = 0.284494
$ws.Cells.Item(19, 4).Value = 8719626932
$ws.Cells.Item(19, 5).Value = 379476851
$ws.Cells.Item(19, 6).Value = -10.99505

$ws.Cells.Item(20, 1).Value = "Avalanche"
$ws.Cells.Item(20, 2).Value = "AVAX"
$ws.Cells.Item(20, 3).Value = 19.73
$ws.Cells.Item(20, 4).Value = 8168511203
$ws.Cells.Item(20, 5).Value = 859886957
$ws.Cells.Item(20, 6).Value = -15.08448

$ws.Cells.Item(21, 1).Value = "USDS"
$ws.Cells.Item(21, 2).Value = "USDS"
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 7931597705
$ws.Cells.Item(21, 5).Value = 87782461
$ws.Cells.Item(21, 6).Value = 0.11705

$ws.Cells.Item(22, 1).Value = "Toncoin"
$ws.Cells.Item(22, 2).Value = "TON"
$ws.Cells.Item(22, 3).Value = 3.07
$ws.Cells.Item(22, 4).Value = 7582678216
$ws.Cells.Item(22, 5).Value = 314288099
$ws.Cells.Item(22, 6).Value = -6.23516

$ws.Cells.Item(23, 1).Value = "Litecoin"
$ws.Cells.Item(23, 2).Value = "LTC"
$ws.Cells.Item(23, 3).Value = 100.53
$ws.Cells.Item(23, 4).Value = 7582361210
$ws.Cells.Item(23, 5).Value = 1667223336
$ws.Cells.Item(23, 6).Value = -12.1806

$ws.Cells.Item(24, 1).Value = "Shiba Inu"
$ws.Cells.Item(24, 2).Value = "SHIB"
$ws.Cells.Item(24, 3).Value = 0.00001271
$ws.Cells.Item(24, 4).Value = 7485464763
$ws.Cells.Item(24, 5).Value = 499213353
$ws.Cells.Item(24, 6).Value = -7.29407

$ws.Cells.Item(25, 1).Value = "Sui"
$ws.Cells.Item(25, 2).Value = "SUI"
$ws.Cells.Item(25, 3).Value = 2.37
$ws.Cells.Item(25, 4).Value = 7482590100
$ws.Cells.Item(25, 5).Value = 1564766675
$ws.Cells.Item(25, 6).Value = -17.29152

$ws.Cells.Item(26, 1).Value = "MANTRA"
$ws.Cells.Item(26, 2).Value = "OM"
$ws.Cells.Item(26, 3).Value = 7.14
$ws.Cells.Item(26, 4).Value = 6969282724
$ws.Cells.Item(26, 5).Value = 228536038
$ws.Cells.Item(26, 6).Value = -3.57664

$ws.Cells.Item(27, 1).Value = "Polkadot"
$ws.Cells.Item(27, 2).Value = "DOT"
$ws.Cells.Item(27, 3).Value = 4.21
$ws.Cells.Item(27, 4).Value = 6396078315
$ws.Cells.Item(27, 5).Value = 408606310
$ws.Cells.Item(27, 6).Value = -13.35549

$ws.Cells.Item(28, 1).Value = "Bitcoin Cash"
$ws.Cells.Item(28, 2).Value = "BCH"
$ws.Cells.Item(28, 3).Value = 305.21
$ws.Cells.Item(28, 4).Value = 6036150222
$ws.Cells.Item(28, 5).Value = 403445875
$ws.Cells.Item(28, 6).Value = -11.58443

$ws.Cells.Item(29, 1).Value = "WETH"
$ws.Cells.Item(29, 2).Value = "WETH"
$ws.Cells.Item(29, 3).Value = 2078.74
$ws.Cells.Item(29, 4).Value = 5997757620
$ws.Cells.Item(29, 5).Value = 1960438435
$ws.Cells.Item(29, 6).Value = -8.698700000000001

$ws.Cells.Item(30, 1).Value = "Ethena USDe"
$ws.Cells.Item(30, 2).Value = "USDE"
$ws.Cells.Item(30, 3).Value = 0.9988899999999999
$ws.Cells.Item(30, 4).Value = 5445353255
$ws.Cells.Item(30, 5).Value = 103661843
$ws.Cells.Item(30, 6).Value = -0.02512

$ws.Cells.Item(31, 1).Value = "Hyperliquid"
$ws.Cells.Item(31, 2).Value = "HYPE"
$ws.Cells.Item(31, 3).Value = 16.29
$ws.Cells.Item(31, 4).Value = 5431024599
$ws.Cells.Item(31, 5).Value = 375742757
$ws.Cells.Item(31, 6).Value = -13.55829

$ws.Cells.Item(32, 1).Value = "Bitget Token"
$ws.Cells.Item(32, 2).Value = "BGB"
$ws.Cells.Item(32, 3).Value = 4.11
$ws.Cells.Item(32, 4).Value = 4925609512
$ws.Cells.Item(32, 5).Value = 344911202
$ws.Cells.Item(32, 6).Value = -7.7355

$ws.Cells.Item(33, 1).Value = "Wrapped eETH"
$ws.Cells.Item(33, 2).Value = "WEETH"
$ws.Cells.Item(33, 3).Value = 2207.92
$ws.Cells.Item(33, 4).Value = 4234814094
$ws.Cells.Item(33, 5).Value = 35391258
$ws.Cells.Item(33, 6).Value = -8.6998

$ws.Cells.Item(34, 1).Value = "WhiteBIT Coin"
$ws.Cells.Item(34, 2).Value = "WBT"
$ws.Cells.Item(34, 3).Value = 29.34
$ws.Cells.Item(34, 4).Value = 4226106083
$ws.Cells.Item(34, 5).Value = 145682879
$ws.Cells.Item(34, 6).Value = -3.24831

$ws.Cells.Item(35, 1).Value = "Uniswap"
$ws.Cells.Item(35, 2).Value = "UNI"
$ws.Cells.Item(35, 3).Value = 6.8
$ws.Cells.Item(35, 4).Value = 4088934005
$ws.Cells.Item(35, 5).Value = 407821404
$ws.Cells.Item(35, 6).Value = -11.54285

$ws.Cells.Item(36, 1).Value = "Monero"
$ws.Cells.Item(36, 2).Value = "XMR"
$ws.Cells.Item(36, 3).Value = 212.38
$ws.Cells.Item(36, 4).Value = 3919602303
$ws.Cells.Item(36, 5).Value = 67227356
$ws.Cells.Item(36, 6).Value = -7.93787

$ws.Cells.Item(37, 1).Value = "Dai"
$ws.Cells.Item(37, 2).Value = "DAI"
$ws.Cells.Item(37, 3).Value = 0.999869
$ws.Cells.Item(37, 4).Value = 3271264529
$ws.Cells.Item(37, 5).Value = 229448209
$ws.Cells.Item(37, 6).Value = 0.00805

$ws.Cells.Item(38, 1).Value = "NEAR Protocol"
$ws.Cells.Item(38, 2).Value = "NEAR"
$ws.Cells.Item(38, 3).Value = 2.69
$ws.Cells.Item(38, 4).Value = 3202373618
$ws.Cells.Item(38, 5).Value = 416906706
$ws.Cells.Item(38, 6).Value = -17.8189

$ws.Cells.Item(39, 1).Value = "Aptos"
$ws.Cells.Item(39, 2).Value = "APT"
$ws.Cells.Item(39, 3).Value = 5.26
$ws.Cells.Item(39, 4).Value = 3098894462
$ws.Cells.Item(39, 5).Value = 413756015
$ws.Cells.Item(39, 6).Value = -15.42943

$ws.Cells.Item(40, 1).Value = "Aave"
$ws.Cells.Item(40, 2).Value = "AAVE"
$ws.Cells.Item(40, 3).Value = 199.64
$ws.Cells.Item(40, 4).Value = 3011983197
$ws.Cells.Item(40, 5).Value = 687611429
$ws.Cells.Item(40, 6).Value = 0.14094

$ws.Cells.Item(41, 1).Value = "sUSDS"
$ws.Cells.Item(41, 2).Value = "SUSDS"
$ws.Cells.Item(41, 3).Value = 1.043
$ws.Cells.Item(41, 4).Value = 3007240928
$ws.Cells.Item(41, 5).Value = 1878851
$ws.Cells.Item(41, 6).Value = 0.26927

$ws.Cells.Item(42, 1).Value = "Ondo"
$ws.Cells.Item(42, 2).Value = "ONDO"
$ws.Cells.Item(42, 3).Value = 0.914899
$ws.Cells.Item(42, 4).Value = 2893183374
$ws.Cells.Item(42, 5).Value = 390923395
$ws.Cells.Item(42, 6).Value = -17.13825

$ws.Cells.Item(43, 1).Value = "Internet Computer"
$ws.Cells.Item(43, 2).Value = "ICP"
$ws.Cells.Item(43, 3).Value = 5.86
$ws.Cells.Item(43, 4).Value = 2821650407
$ws.Cells.Item(43, 5).Value = 132799129
$ws.Cells.Item(43, 6).Value = -12.43504

$ws.Cells.Item(44, 1).Value = "Pepe"
$ws.Cells.Item(44, 2).Value = "PEPE"
$ws.Cells.Item(44, 3).Value = 0.0000067
$ws.Cells.Item(44, 4).Value = 2813164120
$ws.Cells.Item(44, 5).Value = 1176159725
$ws.Cells.Item(44, 6).Value = -15.55409

$ws.Cells.Item(45, 1).Value = "Ethereum Classic"
$ws.Cells.Item(45, 2).Value = "ETC"
$ws.Cells.Item(45, 3).Value = 18.46
$ws.Cells.Item(45, 4).Value = 2783676811
$ws.Cells.Item(45, 5).Value = 176426353
$ws.Cells.Item(45, 6).Value = -7.39281

$ws.Cells.Item(46, 1).Value = "Gate"
$ws.Cells.Item(46, 2).Value = "GT"
$ws.Cells.Item(46, 3).Value = 19.97
$ws.Cells.Item(46, 4).Value = 2513411243
$ws.Cells.Item(46, 5).Value = 24960346
$ws.Cells.Item(46, 6).Value = -7.37125

$ws.Cells.Item(47, 1).Value = "OKB"
$ws.Cells.Item(47, 2).Value = "OKB"
$ws.Cells.Item(47, 3).Value = 41.21
$ws.Cells.Item(47, 4).Value = 2470910297
$ws.Cells.Item(47, 5).Value = 28268635
$ws.Cells.Item(47, 6).Value = -10.7161

$ws.Cells.Item(48, 1).Value = "Official Trump"
$ws.Cells.Item(48, 2).Value = "TRUMP"
$ws.Cells.Item(48, 3).Value = 12.38
$ws.Cells.Item(48, 4).Value = 2470544227
$ws.Cells.Item(48, 5).Value = 1980909099
$ws.Cells.Item(48, 6).Value = -15.52849

$ws.Cells.Item(49, 1).Value = "Coinbase Wrapped BTC"
$ws.Cells.Item(49, 2).Value = "CBBTC"
$ws.Cells.Item(49, 3).Value = 83911
$ws.Cells.Item(49, 4).Value = 2421173266
$ws.Cells.Item(49, 5).Value = 538985057
$ws.Cells.Item(49, 6).Value = -6.89961

$ws.Cells.Item(50, 1).Value = "Mantle"
$ws.Cells.Item(50, 2).Value = "MNT"
$ws.Cells.Item(50, 3).Value = 0.715885
$ws.Cells.Item(50, 4).Value = 2406940619
$ws.Cells.Item(50, 5).Value = 116363262
$ws.Cells.Item(50, 6).Value = -4.5626

$ws.Cells.Item(51, 1).Value = "Tokenize Xchange"
$ws.Cells.Item(51, 2).Value = "TKX"
$ws.Cells.Item(51, 3).Value = 28.86
$ws.Cells.Item(51, 4).Value = 2304540692
$ws.Cells.Item(51, 5).Value = 14262380
$ws.Cells.Item(51, 6).Value = -2.35183

$ws2 = $wb.Worksheets.Item("Analysis")
$ws2.Cells.Item(2, 2).Value = "2025-03-04 22:47:56"
$ws2.Cells.Item(4, 2).NumberFormat = "@"
$ws2.Cells.Item(4, 2).Value = "`$5286.67"
$ws2.Cells.Item(5, 2).Value = "Pi Network (1.24%)"
$ws2.Cells.Item(6, 2).Value = "NEAR Protocol (-17.82%)"

Write-Host "Edit applied"